# Add files via upload
#
# "All structural tables have been made uniform in terms of how valency
# classes, X and Y columns etc. are organized."
#
# Column I is "X" (the valency class of the first argument) and column J
# is "Y" (the valency class of the second argument). A number of rows had
# column I filled in but left column J blank even though the pattern was
# already implied:
#   - I = "TR" ("transitive")         -> I becomes "NOM", J becomes "ACC"
#   - I = "*"  (pattern not assigned) -> J is mirrored to "*" as well
#
# Also reset the sheet view back to the top-left / A1 instead of the
# scrolled-to-K1, M2:M131-selected state it was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slovenian")

# Rows where column I ("X") reads "TR" -> becomes "NOM", column J ("Y")
# becomes "ACC".
$trRows = @(3,5,9,10,12,14,16,17,19,20,21,22,27,28,29,32,34,35,37,39,40,41,42,43,44,45,47,50,51,53,54,56,58,60,61,64,67,69,70,71,72,73,75,76,79,86,87,88,92,93,94,96,97,98,101,103,104,106,107,108,109,111,120,122,125,127)

foreach ($r in $trRows) {
    $cellI = $ws.Cells.Item($r, 9)
    $cellJ = $ws.Cells.Item($r, 10)
    if ($cellI.Text -eq "TR") {
        $cellI.Value = "NOM"
        $cellJ.Value = "ACC"
    }
}

# Rows where column I ("X") reads "*" -> column J ("Y") is filled in to
# match ("*").
$starRows = @(6,66,129)

foreach ($r in $starRows) {
    $cellI = $ws.Cells.Item($r, 9)
    $cellJ = $ws.Cells.Item($r, 10)
    if ($cellI.Text -eq "*") {
        $cellJ.Value = "*"
    }
}

# Reset the view: scroll back to the top-left corner and collapse the
# selection down to a single cell (A1) instead of the M2:M131 block.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()
